# Re-save edit:
#  1. Swap the order of the "Agenda" and "Info" slides (positions 4 and 5).
#  2. Merge the two text runs in the "Caveat" slide's second bullet into a
#     single run (removing the stray run split / smtClean artifact).

$p = $ppt.ActivePresentation

# --- 1. Swap slides 4 ("Agenda") and 5 ("Info") -----------------------------
$agendaSlide = $p.Slides.Item(4)
$agendaSlide.MoveTo(5)

# --- 2. Merge the SSO sentence runs on the "Caveat" slide -------------------
$caveat = $p.Slides.Item(3)
$body = $caveat.Shapes.Item(2)
$tr = $body.TextFrame.TextRange
$ssoPara = $tr.Paragraphs(2, 1)

# Force a genuine text rewrite (identical-text assignment is a no-op) so the
# run split collapses into a single run.
$ssoPara.Text = "PLACEHOLDER"
$ssoPara2 = $tr.Paragraphs(2, 1)
$ssoPara2.Text = "Implementation of an SSO solution will vary depending on your business needs."
